$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: '56 x 91' -> '19 x 78'
$cell = $t.Rows.Item(1).Cells.Item(1)
$cell.Range.Text = "19 x 78" + [char]11 + "  7    8" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "9|    |"

# Row 1, Col 2: '26 x 48' -> '89 x 33'
$cell = $t.Rows.Item(1).Cells.Item(2)
$cell.Range.Text = "89 x 33" + [char]11 + "  3    3" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"

# Row 1, Col 3: '86 x 83' -> '72 x 43'
$cell = $t.Rows.Item(1).Cells.Item(3)
$cell.Range.Text = "72 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "2|    |"

# Row 2, Col 1: '62 x 74' -> '56 x 46'
$cell = $t.Rows.Item(2).Cells.Item(1)
$cell.Range.Text = "56 x 46" + [char]11 + "  4    6" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "6|    |"

# Row 2, Col 2: '44 x 17' -> '57 x 68'
$cell = $t.Rows.Item(2).Cells.Item(2)
$cell.Range.Text = "57 x 68" + [char]11 + "  6    8" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "7|    |"

# Row 2, Col 3: '30 x 98' -> '90 x 29'
$cell = $t.Rows.Item(2).Cells.Item(3)
$cell.Range.Text = "90 x 29" + [char]11 + "  2    9" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "0|    |"

# Row 3, Col 1: '54 x 71' -> '13 x 46'
$cell = $t.Rows.Item(3).Cells.Item(1)
$cell.Range.Text = "13 x 46" + [char]11 + "  4    6" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "3|    |"

# Row 3, Col 2: '73 x 12' -> '13 x 34'
$cell = $t.Rows.Item(3).Cells.Item(2)
$cell.Range.Text = "13 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "3|    |"

# Row 3, Col 3: '64 x 77' -> '30 x 86'
$cell = $t.Rows.Item(3).Cells.Item(3)
$cell.Range.Text = "30 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "0|    |"

# Row 4, Col 1: '83 x 39' -> '99 x 48'
$cell = $t.Rows.Item(4).Cells.Item(1)
$cell.Range.Text = "99 x 48" + [char]11 + "  4    8" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "9|    |"

# Row 4, Col 2: '90 x 98' -> '39 x 51'
$cell = $t.Rows.Item(4).Cells.Item(2)
$cell.Range.Text = "39 x 51" + [char]11 + "  5    1" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "9|    |"

# Row 4, Col 3: '55 x 30' -> '16 x 45'
$cell = $t.Rows.Item(4).Cells.Item(3)
$cell.Range.Text = "16 x 45" + [char]11 + "  4    5" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"

# Row 5, Col 1: '88 x 14' -> '59 x 50'
$cell = $t.Rows.Item(5).Cells.Item(1)
$cell.Range.Text = "59 x 50" + [char]11 + "  5    0" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "9|    |"

# Row 5, Col 2: '95 x 74' -> '65 x 16'
$cell = $t.Rows.Item(5).Cells.Item(2)
$cell.Range.Text = "65 x 16" + [char]11 + "  1    6" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "5|    |"

# Row 5, Col 3: '84 x 35' -> '38 x 31'
$cell = $t.Rows.Item(5).Cells.Item(3)
$cell.Range.Text = "38 x 31" + [char]11 + "  3    1" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "8|    |"

Write-Host "Updated 15 cells"